# Projetos Aula01 solucao ex4
#
# Fill in the "PONTOS" (scores) sheet, column B, rows 3-20, with each
# student's score out of 10. Row 3 (B3) previously held only a stray
# underline character-style with no value, so clear that formatting before
# writing its score. Finish with the selection on B5, matching where the
# author ended up after entering the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PONTOS")
$ws.Activate()

# B3 carried a leftover underline style (no value yet) - drop it before
# typing the score so the cell reads like the rest of the column.
$ws.Range("B3").Font.Underline = $false

$scores = @(10, 10, 5, 10, 10, 0, 0, 10, 10, 10, 10, 10, 10, 10, 10, 0, 10, 10)

for ($i = 0; $i -lt $scores.Length; $i++) {
    $row = 3 + $i
    $ws.Cells.Item($row, 2).Value = $scores[$i]
}

$ws.Range("B5").Select()
